$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hourly refresh of the cryptos list (coinranking.com scrape): update the
# Price (D) and Volume(1h) (E) columns for each coin row. Prices are kept
# as literal text (e.g. "27.074.32", "1.000") rather than locale numbers,
# so each Price cell has its NumberFormat forced to "@" (Text) right before
# the write, then restored to the default "Normal" style afterwards -
# otherwise Excel would silently reinterpret a value like "1.000" as the
# number 1.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.074.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.895.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.24%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5191'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3773'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07218'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8903'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07670'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.898.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.224'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008507'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '

$ws.Range("E18").Value = '  +1.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.132.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.40%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.060'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.126.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.404'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.290'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.957'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.791'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09189'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05046'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.237'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7774'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.973'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.297'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.580'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5630'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01990'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.071'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.973'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.623'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.53%  '

$ws.Range("E45").Value = '  +2.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4827'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.06%  '

# Rows 47-48: coin ranking reshuffled - EnergySwap now appears above PaxDollar
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.56%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.594'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.35%  '
